$d = $word.ActiveDocument
$t = $d.Tables.Item(2)
$cell = $t.Cell(14, 2)
$cr = $cell.Range
$p = $cr.Paragraphs(1)
$pr = $p.Range
$beforeMark = $d.Range($pr.End - 1, $pr.End - 1)

$snippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="259" w:lineRule="auto"/><w:ind w:left="13" w:right="0" w:hanging="11"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Overall code quality</w:t></w:r></w:p></w:body></w:document>'
$beforeMark.InsertXML($snippet)
Write-Output "done"
